$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data row
$ws.Rows.Item(1).Insert()

# New header row (row 1)
$ws.Range("A1").Value = "fname"
$ws.Range("B1").Value = "lname"
$ws.Range("C1").Value = "postcode"

# Existing data moved to row 2; update postcode to text value "E12312"
$ws.Range("C2").Value = "E12312"

# Reset the active selection back to the top-left cell
$ws.Range("A1").Select()
